$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.515473961830139
$ws.Range("B1").Value = 2.007537364959717
$ws.Range("C1").Value = 2.089452981948853
$ws.Range("D1").Value = 1.634312868118286
$ws.Range("E1").Value = 1.455596804618835
